$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 7: update Commit Comment (C7) ---
$ws.Range("C7").Value = "fixed Railroad legend/fixed Area Office misspelling/added link to minute orders"

# --- Row 8: change A8 from text "11/17/" to a real date 11/17/2015, update D8, add C8/E8 ---
# (set D8 before C8 so new shared strings are appended in the same order Excel originally wrote them)
$ws.Range("A8").Value = (Get-Date -Year 2015 -Month 11 -Day 17 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D8").Value = "Added renderer to Area Office layer to mimic transparency on old layer`nFixed RR popup title`nremoved ESRI logo"
$ws.Range("C8").Value = "added area office renderer/fixed RR popup title/ removed ESRI logo"
$ws.Range("E8").Value = "see SPM_TestScrip.docx"

# --- Row 9: new row ---
$ws.Range("A9").Value = (Get-Date -Year 2015 -Month 12 -Day 1 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B9").Value = "cbardash"
$ws.Range("D9").Value = "Fixed ""Clear Overlays"" legend.`nAdded new 2015 MPO Map Service layer and added renderer to handle transparency problem.`nReplaced old service layer for Highway Designations with new Map Service (now includes historic minute orders and has been resymbolized).`nAdded highlight function to Search tab."
$ws.Range("C9").Value = "fixed legend for Clear Overlays/updated MPO layer/updatedHighwayDesignations layer/added highlighter for Search tab"
$ws.Range("E9").Value = "see SPM_TestScrip.docx"
$ws.Range("F9").Value = "YES"

# --- Apply the same formatting as row 7 to the new/changed rows 8 and 9 ---
$ws.Range("A7:F7").Copy() | Out-Null
$ws.Range("A8:F8").PasteSpecial(-4122) | Out-Null
$ws.Range("A9:F9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Row heights ---
$ws.Rows("8").RowHeight = 45
$ws.Rows("9").RowHeight = 90

# --- Update sheet view: scroll so row 5 is at top, select C9 ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 5
$ws.Range("C9").Select() | Out-Null
